$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (copy format from H1, the last existing header cell).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the I and J data columns for rows 2-75.
# Each tuple is: row, I-value, J-value
$data = @(
    @(2, 9, 9),
    @(3, 10, 10),
    @(4, 9, 9),
    @(5, 8, 8),
    @(6, 9, 9),
    @(7, 9, 9),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 8, 8),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 9, 9),
    @(14, 8, 8),
    @(15, 9, 9),
    @(16, 6, 6),
    @(17, 6, 6),
    @(18, 7, 7),
    @(19, 7, 7),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 6, 7),
    @(23, 7, 7),
    @(24, 8, 8),
    @(25, 7, 7),
    @(26, 8, 8),
    @(27, 6, 6),
    @(28, 6, 7),
    @(29, 6, 6),
    @(30, 6, 6),
    @(31, 7, 7),
    @(32, 8, 8),
    @(33, 7, 7),
    @(34, 8, 8),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 7, 7),
    @(38, 7, 7),
    @(39, 7, 7),
    @(40, 6, 7),
    @(41, 6, 7),
    @(42, 7, 8),
    @(43, 6, 6),
    @(44, 7, 7),
    @(45, 5, 5),
    @(46, 6, 6),
    @(47, 9, 9),
    @(48, 6, 7),
    @(49, 7, 7),
    @(50, 8, 8),
    @(51, 7, 7),
    @(52, 5, 5),
    @(53, 2, 4),
    @(54, 5, 6),
    @(55, 9, 9),
    @(56, 6, 6),
    @(57, 8, 8),
    @(58, 7, 7),
    @(59, 7, 7),
    @(60, 8, 9),
    @(61, 6, 7),
    @(62, 7, 7),
    @(63, 6, 6),
    @(64, 8, 8),
    @(65, 7, 7),
    @(66, 7, 7),
    @(67, 6, 6),
    @(68, 7, 7),
    @(69, 8, 9),
    @(70, 7, 7),
    @(71, 7, 7),
    @(72, 1, 2),
    @(73, 8, 8),
    @(74, 6, 6),
    @(75, 7, 7)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Host "Added I0/IF columns: $($data.Count) data rows updated."
